$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update the "Student ID" (column A) values for the existing rows 2-14.
#    A leading apostrophe forces these numeric-looking values to be stored
#    as text, matching the original inline-string typing of the column.
# ---------------------------------------------------------------------------
$idUpdates = @{
    2  = "191061"
    3  = "191502"
    4  = "191480"
    5  = "191375"
    6  = "181013"
    7  = "211175"
    8  = "201218"
    9  = "201297"
    10 = "201253"
    11 = "200727"
    12 = "200708"
    13 = "200943"
    14 = "201343"
}

foreach ($row in $idUpdates.Keys) {
    $ws.Cells.Item($row, 1).Value = "'" + $idUpdates[$row]
}

# ---------------------------------------------------------------------------
# 2) Append five new log rows (15-19), all sharing the same Subject/Date/
#    Time/Type/User values, alternating the same banded row styling used by
#    the existing table (row 13's style for odd rows, row 14's for even).
# ---------------------------------------------------------------------------
$newRows = @(
    @{ Row = 15; Id = "200405"; FormatSrc = 13 }
    @{ Row = 16; Id = "201682"; FormatSrc = 14 }
    @{ Row = 17; Id = "200858"; FormatSrc = 13 }
    @{ Row = 18; Id = "200938"; FormatSrc = 14 }
    @{ Row = 19; Id = "201065"; FormatSrc = 13 }
)

foreach ($entry in $newRows) {
    $r = $entry.Row
    $srcRow = $entry.FormatSrc

    # Clone the banded-row formatting from the matching existing row.
    $ws.Range("A$srcRow`:F$srcRow").Copy()
    $ws.Range("A$r`:F$r").PasteSpecial(-4122)

    # Populate the row's values.
    $ws.Cells.Item($r, 1).Value = "'" + $entry.Id
    $ws.Cells.Item($r, 2).Value = "general surgery"
    $ws.Cells.Item($r, 3).Value = "13/10/2025"
    $ws.Cells.Item($r, 4).Value = "10:30:00"
    $ws.Cells.Item($r, 5).Value = "Excuse"
    $ws.Cells.Item($r, 6).Value = "System"
}
